$wb = $excel.ActiveWorkbook

# --- Sheet "About" ---
$ws1 = $wb.Worksheets.Item("About")

# Clear old content that is being relocated/replaced (old "Note:" at A5 and
# old "Only one quality tier..." note at A6) before rebuilding rows 4-17.
$ws1.Range("A5").Clear()
$ws1.Range("A6").Clear()

# Row 3: Source info
$ws1.Range("B3").Value = "Datta, S. and Filippini, M."

# Row 4: Year (left aligned number)
$ws1.Range("B4").Value = 2012
$ws1.Range("B4").HorizontalAlignment = -4131

# Row 5: Title of paper
$ws1.Range("B5").Value = "Analysing the Impact of Energy Star Rebate Policies in the U.S."

# Row 6: Hyperlink to source PDF
$ws1.Range("B6").Value = "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf"
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf")

# Row 7: Page reference
$ws1.Range("B7").Value = "Page 12, first paragrah"

# Row 9: Note: header
$ws1.Range("A9").Value = "Note:"
$ws1.Range("A9").Font.Bold = $true

# Row 10: variable description
$ws1.Range("A10").Value = "This variable measures how a rebate program influences market shares of rebate-qualifying and non-qualifying components."

# Row 11: analysis-based note
$ws1.Range("A11").Value = "Analysis based on air conditioners, clothes washers, dishwashers, and refrigerators because data not available for other product types (p. 8)"

# Row 13-14: percentage-point clarifications
$ws1.Range("A13").Value = "Note that these figures are expressed as percentage points, not as a percentage growth (or shrinkage) relative to a BAU quantity."
$ws1.Range("A14").Value = "That is, this is a percentage of all sales, not a percentage of the sales of a particular quality level."

# Row 16-17: quality-level note
$ws1.Range("A16").Value = "Changes across all quality levels must sum to zero.  The U.S. dataset only has two quality levels, so the"
$ws1.Range("A17").Value = 'decrease in the "standard-compliant" share must equal the increase in the "rebate-qualifying" share.'

# --- Sheet "MSCdtRPbQL" ---
$ws2 = $wb.Worksheets.Item("MSCdtRPbQL")

$ws2.Range("A2").Value = "Change in Perc Share (dimensionless)"
$ws2.Range("A2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 30

$ws2.Range("C2").Value = 0.074
$ws2.Range("B2").Formula = "=-C2"

# Leave the selection cursor on A2 in MSCdtRPbQL, but keep "About" as the
# active/selected tab (matches original workbook's tabSelected on "About").
$ws2.Range("A2").Select()
$ws1.Activate()
